# Gradle浅析.pptx — 插件调整 by qiyu 2021-12-15
#
# 1) Fix a typo on the "增量式构建（Up-to-date checks）" slide:
#      检查输入输入输出是否变化...  ->  检查输入输出是否变化...
# 2) Append a new "标题和内容" slide titled “依赖管理” at the end of the deck.

$p = $ppt.ActivePresentation

# --- 1) Typo fix on the existing slide ------------------------------------
$typoSlide = $p.Slides.Item(18)
$typoShape = $typoSlide.Shapes.Item(2)
$typoRange = $typoShape.TextFrame.TextRange

$oldText = "检查输入输入输出是否变化（文件大小和最后更新时间）来判断是否执行增量式构建。"
$newText = "检查输入输出是否变化（文件大小和最后更新时间）来判断是否执行增量式构建。"

$fullText = $typoRange.Text
$startPos = $fullText.IndexOf($oldText)
if ($startPos -ge 0) {
    $run = $typoRange.Characters($startPos + 1, $oldText.Length)
    $run.Text = $newText
}

# --- 2) New "依赖管理" slide -------------------------------------------------
# Slides 2-19 all use the master's 2nd custom layout ("标题和内容").
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $titleAndContent)

$newSlide.Shapes.Item(1).Name = "标题 1"
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "依赖管理"

$newSlide.Shapes.Item(2).Name = "内容占位符 2"
